# Scheduled runner update: refresh market-price / profit columns (H:N)
# on the leve tables across several job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3833.7778
$ws.Range("I62").Value = 4534
$ws.Range("K62").Value = 4534
$ws.Range("M62").Value = -3910

$ws.Range("H65").Value = 3833.7778
$ws.Range("I65").Value = 4534
$ws.Range("K65").Value = 22670
$ws.Range("M65").Value = -19550

$ws.Range("H98").Value = 4561.846
$ws.Range("I98").Value = 2441
$ws.Range("J98").Value = 11631.333
$ws.Range("K98").Value = 2441
$ws.Range("L98").Value = 11631.333
$ws.Range("M98").Value = -943
$ws.Range("N98").Value = -14627.333

$ws.Range("H122").Value = 4561.846
$ws.Range("I122").Value = 2441
$ws.Range("J122").Value = 11631.333
$ws.Range("K122").Value = 7323
$ws.Range("L122").Value = 34893.999
$ws.Range("M122").Value = -4873
$ws.Range("N122").Value = -39793.999

$ws.Range("H127").Value = 1078028.8
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 1078028.8
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 3234086.4
$ws.Range("M127").ClearContents() | Out-Null
$ws.Range("N127").Value = -3244006.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10257.164
$ws.Range("I32").Value = 11228.949
$ws.Range("J32").Value = 6161.7856
$ws.Range("K32").Value = 11228.949
$ws.Range("L32").Value = 6161.7856
$ws.Range("M32").Value = -10941.949
$ws.Range("N32").Value = -6735.7856

$ws.Range("H74").Value = 1046.5349
$ws.Range("I74").Value = 979.17145
$ws.Range("J74").Value = 1341.25
$ws.Range("K74").Value = 979.17145
$ws.Range("L74").Value = 1341.25
$ws.Range("M74").Value = -105.17145
$ws.Range("N74").Value = -3089.25

$ws.Range("H77").Value = 1046.5349
$ws.Range("I77").Value = 979.17145
$ws.Range("J77").Value = 1341.25
$ws.Range("K77").Value = 4895.85725
$ws.Range("L77").Value = 6706.25
$ws.Range("M77").Value = -527.85725
$ws.Range("N77").Value = -15442.25

$ws.Range("H110").Value = 1987.25
$ws.Range("I110").Value = 1987.25
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1987.25
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 57.75
$ws.Range("N110").ClearContents() | Out-Null

$ws.Range("H122").Value = 9256.549999999999
$ws.Range("I122").Value = 9101.723
$ws.Range("J122").Value = 10650
$ws.Range("K122").Value = 27305.169
$ws.Range("L122").Value = 31950
$ws.Range("M122").Value = -24855.169
$ws.Range("N122").Value = -36850

$ws.Range("H132").Value = 22730312
$ws.Range("I132").Value = 62502930
$ws.Range("J132").Value = 3103.7144
$ws.Range("K132").Value = 187508790
$ws.Range("L132").Value = 9311.143199999999
$ws.Range("M132").Value = -187506260
$ws.Range("N132").Value = -14371.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 27780550
$ws.Range("I86").Value = 2492.8572
$ws.Range("J86").Value = 125003750
$ws.Range("K86").Value = 2492.8572
$ws.Range("L86").Value = 125003750
$ws.Range("M86").Value = -1369.8572
$ws.Range("N86").Value = -125005996

$ws.Range("H89").Value = 27780550
$ws.Range("I89").Value = 2492.8572
$ws.Range("J89").Value = 125003750
$ws.Range("K89").Value = 12464.286
$ws.Range("L89").Value = 625018750
$ws.Range("M89").Value = -6848.286
$ws.Range("N89").Value = -625029982

$ws.Range("H94").Value = 879.6667
$ws.Range("I94").Value = 761.5833
$ws.Range("J94").Value = 1115.8334
$ws.Range("K94").Value = 761.5833
$ws.Range("L94").Value = 1115.8334
$ws.Range("M94").Value = -310.5833
$ws.Range("N94").Value = -2017.8334

$ws.Range("H105").Value = 3378.45
$ws.Range("I105").Value = 1585.4445
$ws.Range("J105").Value = 4845.4546
$ws.Range("K105").Value = 1585.4445
$ws.Range("L105").Value = 4845.4546
$ws.Range("M105").Value = 161.5554999999999
$ws.Range("N105").Value = -8339.454600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 796159.1
$ws.Range("I134").Value = 2634.8333
$ws.Range("J134").Value = 3970256.2
$ws.Range("K134").Value = 7904.499899999999
$ws.Range("L134").Value = 11910768.6
$ws.Range("M134").Value = -5369.499899999999
$ws.Range("N134").Value = -11915838.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1836.3334
$ws.Range("I122").Value = 1474.6
$ws.Range("J122").Value = 2559.8
$ws.Range("K122").Value = 4423.799999999999
$ws.Range("L122").Value = 7679.400000000001
$ws.Range("M122").Value = -1973.799999999999
$ws.Range("N122").Value = -12579.4

$ws.Range("H132").Value = 4415.2856
$ws.Range("I132").Value = 4273.0713
$ws.Range("J132").Value = 4699.7144
$ws.Range("K132").Value = 12819.2139
$ws.Range("L132").Value = 14099.1432
$ws.Range("M132").Value = -10289.2139
$ws.Range("N132").Value = -19159.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1124.0264
$ws.Range("I22").Value = 535.0833
$ws.Range("J22").Value = 1395.8462
$ws.Range("K22").Value = 535.0833
$ws.Range("L22").Value = 1395.8462
$ws.Range("M22").Value = -240.0833
$ws.Range("N22").Value = -1985.8462

$ws.Range("H27").Value = 1124.0264
$ws.Range("I27").Value = 535.0833
$ws.Range("J27").Value = 1395.8462
$ws.Range("K27").Value = 535.0833
$ws.Range("L27").Value = 1395.8462
$ws.Range("M27").Value = -428.0833
$ws.Range("N27").Value = -1609.8462

$ws.Range("H40").Value = 6035.92
$ws.Range("I40").Value = 6999.778
$ws.Range("J40").Value = 5493.75
$ws.Range("K40").Value = 6999.778
$ws.Range("L40").Value = 5493.75
$ws.Range("M40").Value = -6863.778
$ws.Range("N40").Value = -5765.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2200.6667
$ws.Range("I107").Value = 2520.8
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 7562.400000000001
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = -5642.400000000001
$ws.Range("N107").Value = -5640

$ws.Range("H122").Value = 2718.5454
$ws.Range("I122").Value = 2840.4
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 8521.200000000001
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -6071.200000000001
$ws.Range("N122").Value = -9400

$ws.Range("H132").Value = 1138.2559
$ws.Range("I132").Value = 951.0769
$ws.Range("J132").Value = 2963.25
$ws.Range("K132").Value = 2853.2307
$ws.Range("L132").Value = 8889.75
$ws.Range("M132").Value = -323.2307000000001
$ws.Range("N132").Value = -13949.75

$ws.Range("H136").Value = 2045.3636
$ws.Range("I136").Value = 2055.4443
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 6166.3329
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -3616.3329
$ws.Range("N136").Value = -11100

$ws.Range("H137").Value = 20905
$ws.Range("I137").Value = 20000
$ws.Range("J137").Value = 21357.5
$ws.Range("K137").Value = 20000
$ws.Range("L137").Value = 21357.5
$ws.Range("M137").Value = -14900
$ws.Range("N137").Value = -31557.5
